# Append 11 new daily JPY/USD FX rows (2023-10-30 .. 2023-11-13) to Sheet1,
# extending the data table from row 855 to row 866.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, extend formatting (style) of the last existing data row (855) down
# through the new rows so the new date cells in column A pick up the same
# date-formatted style (s="2") as the rest of the table, and the numeric
# columns B:F stay unstyled, exactly like the existing rows.
$ws.Range("A855:F855").Copy()
$ws.Range("A856:F866").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row data: date (Excel serial), open, high, low, close, volume
$newRows = @(
    @(856, 45229, 0.006679, 0.006719, 0.006674, 0.006706, 0),
    @(857, 45230, 0.006705, 0.006708, 0.00659,  0.00659,  0),
    @(858, 45231, 0.006591, 0.006636, 0.006591, 0.006623, 0),
    @(859, 45232, 0.006623, 0.006673, 0.006623, 0.006644, 0),
    @(860, 45233, 0.006644, 0.0067,   0.006644, 0.006692, 0),
    @(861, 45236, 0.006692, 0.006699, 0.006661, 0.006661, 0),
    @(862, 45237, 0.006661, 0.006669, 0.006637, 0.006648, 0),
    @(863, 45238, 0.006648, 0.006651, 0.006621, 0.006622, 0),
    @(864, 45239, 0.00662,  0.006632, 0.006606, 0.006607, 0),
    @(865, 45240, 0.006607, 0.006612, 0.006596, 0.006596, 0),
    @(866, 45243, 0.006597, 0.006604, 0.006591, 0.006591, 0)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}

Write-Host "Appended $($newRows.Count) rows (856-866) to $($ws.Name)."
